$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "29.241.73"
$ws.Range("E2").Value = "  -0.53%  "

$ws.Range("D3").Value = "1.860.68"
$ws.Range("E3").Value = "  -0.96%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.50%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6984"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.54%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.0000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07828"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.26%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3122"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.05"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.74%  "

$ws.Range("E11").Value = "  -4.15%  "

$ws.Range("D12").Value = "1.870.65"
$ws.Range("E12").Value = "  -0.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.133"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.98%  "

$ws.Range("E14").Value = "  -3.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6938"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.91%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.597"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.24%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008523"
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").Value = "29.261.56"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "248.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.52%  "

$ws.Range("D20").Value = "2.114.15"
$ws.Range("E20").Value = "  -1.49%  "

$ws.Range("E21").Value = "  -3.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.0000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.573"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1540"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.908"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.72%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.577"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.81%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.285"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.90%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.237"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.208"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.83%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05233"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.86%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7621"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.52%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.880"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.35%  "

$ws.Range("E36").Value = "  +0.11%  "

$ws.Range("E37").Value = "  -0.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01857"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.80%  "

$ws.Range("D39").Value = "1.242.10"
$ws.Range("E39").Value = "  -2.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.742"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.64%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9017"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.26%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "110.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.901"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.73%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "69.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.94%  "

$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9996"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.15%  "

$ws.Range("D46").Value = "2.011.45"
$ws.Range("E46").Value = "  -1.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000125"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.563"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5185"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.39%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.771"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4263"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.86%  "

